$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# "added BKMP differentiated analysis":
# flag the Beine / Gesäß / Arme columns (L, M, N) in a new row 12 for
# follow-up review.
$ws.Range("L12").Value = "revidieren"
$ws.Range("M12").Value = "revidieren"
$ws.Range("N12").Value = "revidieren"

# Leave the cursor where the author left it when saving.
[void]$ws.Range("G9").Select()
